$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17451536655426
$ws.Range("B1").Value = 2.180500507354736
$ws.Range("C1").Value = 4.433859825134277
$ws.Range("D1").Value = 2.711233854293823
$ws.Range("E1").Value = 1.22624671459198
